$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.047735691070557
$ws.Range("B1").Value = 2.561015367507935
$ws.Range("C1").Value = 2.129439115524292
$ws.Range("D1").Value = 2.008375644683838
$ws.Range("E1").Value = 1.749341368675232
